$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Existing data: A1 = "Scenarios" header, A2 = 1
# Extend the numeric series down to A28 = 27
for ($i = 2; $i -le 27; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $i
}
